$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6543.656484601369
$ws.Range("D3").Value = 270.6623123289115

$ws.Range("B4").Value = 3695.750328977708
$ws.Range("D4").Value = 234.238372227658

$ws.Range("B5").Value = 563.0069999999998

$ws.Range("B6").Value = 3754.008000000001
$ws.Range("D6").Value = 200.001

$ws.Range("B7").Value = 5621.164000000004
$ws.Range("D7").Value = 280.001

$ws.Range("B8").Value = 8845.182250000013
$ws.Range("D8").Value = 294.002

$ws.Range("B9").Value = 28105.07300000001
$ws.Range("D9").Value = 2075

$ws.Range("F10").Value = 3439345848.305001

$ws.Range("G11").Value = 0.7109339862006124

$ws.Range("F12").Value = 142260111.36
$ws.Range("G12").Value = 0.04136254905278265

$ws.Range("G13").Value = 0.2477034647466048
